# Log_of_all_Blogs.xlsx - add Post 57 ("Banker's Algorithm Question | Operating
# System - Mo4 P06") as a new row at the bottom of the Table2 listing on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$lo = $ws.ListObjects.Item(1)

# Duplicate the last data row (row 66) down into the new row 67 so that the
# new row inherits the exact same cell formatting (number format / hyperlink
# style) used throughout the table, then grow the table definition to cover
# the freshly inserted row.
$ws.Rows.Item(66).Copy()
$ws.Rows.Item(67).Insert(-4121)
$excel.CutCopyMode = $false
$lo.Resize($ws.Range("B10:F67"))

# Fill in the values for Post 57. Columns are written F, E, C (then B, D) to
# reproduce the original author's shared-string insertion order.
$ws.Range("B67").Value2 = 57
$ws.Range("F67").Value2 = "https://dev.to/rahulmishra05/banker-s-algorithm-question-operating-system-m04-p06-3lgj"
$ws.Range("E67").Value2 = "https://programmingport.hashnode.dev/bankers-algorithm-question-or-operating-system-m04-p06"
$ws.Range("C67").Value2 = "Banker's Algorithm Question | Operating System - Mo4 P06"
$ws.Range("D67").Value2 = 44177

# Update the view so the newly added row is scrolled into view / selected,
# matching the author's saved cursor position.
$win = $excel.ActiveWindow
$win.ScrollRow = 45
$win.ScrollColumn = 1
$ws.Range("D67").Select()
